# Update NATMI LR-pairs data (Ptn-Ptprz1) with newly recomputed TPM-based
# values. Columns A-D (Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster) are unchanged; columns E:T (all numeric metrics) are
# refreshed with the values produced by the updated scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = [double]"2"
$ws.Cells.Item(2, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2, 7).Value = [double]"0.07024999999999999"
$ws.Cells.Item(2, 8).Value = [double]"0.21075"
$ws.Cells.Item(2, 9).Value = [double]"0.005236595731231519"
$ws.Cells.Item(2, 10).Value = [double]"0.005236595731231519"
$ws.Cells.Item(2, 11).Value = [double]"2"
$ws.Cells.Item(2, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2, 13).Value = [double]"0.01848533333333334"
$ws.Cells.Item(2, 14).Value = [double]"0.05545600000000001"
$ws.Cells.Item(2, 15).Value = [double]"0.001625201930372746"
$ws.Cells.Item(2, 16).Value = [double]"0.001625201930372746"
$ws.Cells.Item(2, 17).Value = [double]"0.001298594666666667"
$ws.Cells.Item(2, 18).Value = [double]"0.011687352"
$ws.Cells.Item(2, 19).Value = [double]"8.510525490979146E-06"
$ws.Cells.Item(2, 20).Value = [double]"8.510525490979144E-06"
$ws.Cells.Item(3, 5).Value = [double]"2"
$ws.Cells.Item(3, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(3, 7).Value = [double]"0.07024999999999999"
$ws.Cells.Item(3, 8).Value = [double]"0.21075"
$ws.Cells.Item(3, 9).Value = [double]"0.005236595731231519"
$ws.Cells.Item(3, 10).Value = [double]"0.005236595731231519"
$ws.Cells.Item(3, 15).Value = [double]"0.002698334581238102"
$ws.Cells.Item(3, 16).Value = [double]"0.002698334581238102"
$ws.Cells.Item(3, 17).Value = [double]"0.002156066166666666"
$ws.Cells.Item(3, 18).Value = [double]"0.0194045955"
$ws.Cells.Item(3, 19).Value = [double]"1.413008734954583E-05"
$ws.Cells.Item(3, 20).Value = [double]"1.413008734954583E-05"
$ws.Cells.Item(4, 5).Value = [double]"2"
$ws.Cells.Item(4, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(4, 7).Value = [double]"0.07024999999999999"
$ws.Cells.Item(4, 8).Value = [double]"0.21075"
$ws.Cells.Item(4, 9).Value = [double]"0.005236595731231519"
$ws.Cells.Item(4, 10).Value = [double]"0.005236595731231519"
$ws.Cells.Item(4, 13).Value = [double]"11.32499966666667"
$ws.Cells.Item(4, 14).Value = [double]"33.974999"
$ws.Cells.Item(4, 15).Value = [double]"0.9956764634883892"
$ws.Cells.Item(4, 16).Value = [double]"0.995676463488389"
$ws.Cells.Item(4, 17).Value = [double]"0.7955812265833334"
$ws.Cells.Item(4, 18).Value = [double]"7.16023103925"
$ws.Cells.Item(4, 19).Value = [double]"0.005213955118390994"
$ws.Cells.Item(4, 20).Value = [double]"0.005213955118390993"
$ws.Cells.Item(5, 9).Value = [double]"0.4287876899474159"
$ws.Cells.Item(5, 10).Value = [double]"0.4287876899474159"
$ws.Cells.Item(5, 11).Value = [double]"2"
$ws.Cells.Item(5, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(5, 13).Value = [double]"0.01848533333333334"
$ws.Cells.Item(5, 14).Value = [double]"0.05545600000000001"
$ws.Cells.Item(5, 15).Value = [double]"0.001625201930372746"
$ws.Cells.Item(5, 16).Value = [double]"0.001625201930372746"
$ws.Cells.Item(5, 17).Value = [double]"0.1063327084764445"
$ws.Cells.Item(5, 18).Value = [double]"0.9569943762880001"
$ws.Cells.Item(5, 19).Value = [double]"0.0006968665814226108"
$ws.Cells.Item(5, 20).Value = [double]"0.0006968665814226107"
$ws.Cells.Item(6, 9).Value = [double]"0.4287876899474159"
$ws.Cells.Item(6, 10).Value = [double]"0.4287876899474159"
$ws.Cells.Item(6, 15).Value = [double]"0.002698334581238102"
$ws.Cells.Item(6, 16).Value = [double]"0.002698334581238102"
$ws.Cells.Item(6, 19).Value = [double]"0.001157012651794314"
$ws.Cells.Item(6, 20).Value = [double]"0.001157012651794313"
$ws.Cells.Item(7, 9).Value = [double]"0.4287876899474159"
$ws.Cells.Item(7, 10).Value = [double]"0.4287876899474159"
$ws.Cells.Item(7, 13).Value = [double]"11.32499966666667"
$ws.Cells.Item(7, 14).Value = [double]"33.974999"
$ws.Cells.Item(7, 15).Value = [double]"0.9956764634883892"
$ws.Cells.Item(7, 16).Value = [double]"0.995676463488389"
$ws.Cells.Item(7, 17).Value = [double]"65.14450490757524"
$ws.Cells.Item(7, 18).Value = [double]"586.3005441681771"
$ws.Cells.Item(7, 19).Value = [double]"0.426933810714199"
$ws.Cells.Item(7, 20).Value = [double]"0.4269338107141989"
$ws.Cells.Item(8, 7).Value = [double]"7.592679666666666"
$ws.Cells.Item(8, 8).Value = [double]"22.778039"
$ws.Cells.Item(8, 9).Value = [double]"0.5659757143213526"
$ws.Cells.Item(8, 10).Value = [double]"0.5659757143213525"
$ws.Cells.Item(8, 11).Value = [double]"2"
$ws.Cells.Item(8, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(8, 13).Value = [double]"0.01848533333333334"
$ws.Cells.Item(8, 14).Value = [double]"0.05545600000000001"
$ws.Cells.Item(8, 15).Value = [double]"0.001625201930372746"
$ws.Cells.Item(8, 16).Value = [double]"0.001625201930372746"
$ws.Cells.Item(8, 17).Value = [double]"0.1403532145315556"
$ws.Cells.Item(8, 18).Value = [double]"1.263178930784"
$ws.Cells.Item(8, 19).Value = [double]"0.000919824823459156"
$ws.Cells.Item(8, 20).Value = [double]"0.0009198248234591558"
$ws.Cells.Item(9, 7).Value = [double]"7.592679666666666"
$ws.Cells.Item(9, 8).Value = [double]"22.778039"
$ws.Cells.Item(9, 9).Value = [double]"0.5659757143213526"
$ws.Cells.Item(9, 10).Value = [double]"0.5659757143213525"
$ws.Cells.Item(9, 15).Value = [double]"0.002698334581238102"
$ws.Cells.Item(9, 16).Value = [double]"0.002698334581238102"
$ws.Cells.Item(9, 17).Value = [double]"0.2330294625428889"
$ws.Cells.Item(9, 18).Value = [double]"2.097265162886"
$ws.Cells.Item(9, 19).Value = [double]"0.001527191842094243"
$ws.Cells.Item(9, 20).Value = [double]"0.001527191842094242"
$ws.Cells.Item(10, 7).Value = [double]"7.592679666666666"
$ws.Cells.Item(10, 8).Value = [double]"22.778039"
$ws.Cells.Item(10, 9).Value = [double]"0.5659757143213526"
$ws.Cells.Item(10, 10).Value = [double]"0.5659757143213525"
$ws.Cells.Item(10, 13).Value = [double]"11.32499966666667"
$ws.Cells.Item(10, 14).Value = [double]"33.974999"
$ws.Cells.Item(10, 15).Value = [double]"0.9956764634883892"
$ws.Cells.Item(10, 16).Value = [double]"0.995676463488389"
$ws.Cells.Item(10, 17).Value = [double]"85.98709469410679"
$ws.Cells.Item(10, 18).Value = [double]"773.8838522469611"
$ws.Cells.Item(10, 19).Value = [double]"0.5635286976557993"
$ws.Cells.Item(10, 20).Value = [double]"0.563528697655799"

Write-Host "Updated LR-pair metrics for rows 2-10 (columns E:T)."
